$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.201.53'
$ws.Range('E2').Value = '  -1.27%  '
$ws.Range('D3').Value = '3.437.84'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''582.92'
$ws.Range('E5').Value = '  -1.68%  '
$ws.Range('D6').Value = '''173.64'
$ws.Range('E6').Value = '  -3.53%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '''0.592'
$ws.Range('E8').Value = '  -3.25%  '
$ws.Range('D9').Value = '3.438.32'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('D10').Value = '''0.130'
$ws.Range('E10').Value = '  -6.49%  '
$ws.Range('E11').Value = '  -1.14%  '
$ws.Range('D12').Value = '''0.410'
$ws.Range('E12').Value = '  -4.78%  '
$ws.Range('D13').Value = '4.038.60'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').Value = '''29.87'
$ws.Range('E15').Value = '  -6.44%  '
$ws.Range('D16').Value = '66.231.54'
$ws.Range('E16').Value = '  -1.18%  '
$ws.Range('E17').Value = '  -3.58%  '
$ws.Range('D18').Value = '3.448.71'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').Value = '''5.88'
$ws.Range('E19').Value = '  -5.16%  '
$ws.Range('D20').Value = '''13.71'
$ws.Range('E20').Value = '  -3.26%  '
$ws.Range('D21').Value = '''373.46'
$ws.Range('E21').Value = '  -4.80%  '
$ws.Range('D22').Value = '''7.73'
$ws.Range('E22').Value = '  -2.37%  '
$ws.Range('D23').Value = '''0.996'
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('D24').Value = '''71.84'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').Value = '''5.70'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('D26').Value = '''0.530'
$ws.Range('E26').Value = '  -1.53%  '
$ws.Range('E27').Value = '  -2.25%  '
$ws.Range('D28').Value = '''9.60'
$ws.Range('E28').Value = '  -7.16%  '
$ws.Range('E29').Value = '  +1.22%  '
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('D31').Value = '''23.95'
$ws.Range('E31').Value = '  +1.76%  '
$ws.Range('D32').Value = '''5.77'
$ws.Range('E32').Value = '  -5.67%  '
$ws.Range('E33').Value = '  -3.01%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '''1.30'
$ws.Range('E35').Value = '  -7.04%  '
$ws.Range('D36').Value = '''7.05'
$ws.Range('E36').Value = '  -3.66%  '
$ws.Range('D37').Value = '''1.54'
$ws.Range('E37').Value = '  -2.32%  '
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('D39').Value = '''29.29'
$ws.Range('E39').Value = '  +12.57%  '
$ws.Range('D40').Value = '''0.885'
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('E41').Value = '  -4.62%  '
$ws.Range('E42').Value = '  -5.81%  '
$ws.Range('D43').Value = '2.731.74'
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('D44').Value = '''4.41'
$ws.Range('E44').Value = '  -5.08%  '
$ws.Range('D45').Value = '''6.28'
$ws.Range('E45').Value = '  -6.47%  '
$ws.Range('D46').Value = '''0.0681'
$ws.Range('E46').Value = '  -5.21%  '
$ws.Range('D47').Value = '''40.25'
$ws.Range('E47').Value = '  -2.49%  '
$ws.Range('D48').Value = '''0.0289'
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('D49').Value = '''23.93'
$ws.Range('E49').Value = '  -8.60%  '
$ws.Range('D50').Value = '''305.75'
$ws.Range('E50').Value = '  -5.76%  '
$ws.Range('D51').Value = '''0.813'
$ws.Range('E51').Value = '  -3.07%  '
